$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price cells whose new values would otherwise be
# auto-converted to numbers by Excel (these are plain text cells in the sheet).
$textCells = @("D5", "D6", "D7", "D8", "D9", "D11", "D14", "D16", "D18", "D20", "D21", "D23", "D25", "D26", "D27", "D29", "D32", "D33", "D37", "D38", "D41", "D42", "D43", "D44", "D46", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '27.461.58'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '1.644.87'
$ws.Range("E3").Value = '  -1.30%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '212.46'
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").Value = '0.529'
$ws.Range("E6").Value = '  +3.42%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '23.37'
$ws.Range("E8").Value = '  -2.96%  '
$ws.Range("D9").Value = '0.257'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("D11").Value = '0.0891'
$ws.Range("E11").Value = '  +1.21%  '
$ws.Range("D12").Value = '1.876.78'
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("D13").Value = '1.626.09'
$ws.Range("E13").Value = '  -2.48%  '
$ws.Range("D14").Value = '0.579'
$ws.Range("E14").Value = '  +2.25%  '
$ws.Range("E15").Value = '  -2.40%  '
$ws.Range("D16").Value = '64.45'
$ws.Range("D17").Value = '27.422.13'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").Value = '230.38'
$ws.Range("E18").Value = '  -4.97%  '
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '7.58'
$ws.Range("E20").Value = '  -1.35%  '
$ws.Range("D21").Value = '1.00'
$ws.Range("E21").Value = '  +0.00%  '
$ws.Range("E22").Value = '  -3.82%  '
$ws.Range("D23").Value = '9.72'
$ws.Range("E23").Value = '  +3.82%  '
$ws.Range("E24").Value = '  -1.29%  '
$ws.Range("D25").Value = '147.64'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").Value = '7.04'
$ws.Range("E26").Value = '  -2.82%  '
$ws.Range("D27").Value = '0.113'
$ws.Range("E27").Value = '  +1.41%  '
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '15.63'
$ws.Range("E30").Value = '  -3.88%  '
$ws.Range("E31").Value = '  -3.73%  '
$ws.Range("D32").Value = '3.29'
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("D33").Value = '3.17'
$ws.Range("E33").Value = '  +1.24%  '
$ws.Range("D34").Value = '1.422.51'
$ws.Range("E34").Value = '  -2.60%  '
$ws.Range("E35").Value = '  +0.43%  '
$ws.Range("E36").Value = '  +0.40%  '
$ws.Range("D37").Value = '0.567'
$ws.Range("E37").Value = '  -1.88%  '
$ws.Range("D38").Value = '0.885'
$ws.Range("E38").Value = '  -4.63%  '
$ws.Range("E39").Value = '  -3.84%  '
$ws.Range("E40").Value = '  -0.84%  '
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '0.823'
$ws.Range("E42").Value = '  +4.00%  '
$ws.Range("D43").Value = '2.46'
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = '5.51'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("E45").Value = '  +0.80%  '
$ws.Range("D46").Value = '64.70'
$ws.Range("E46").Value = '  -7.13%  '
$ws.Range("D47").Value = '1.786.33'
$ws.Range("E47").Value = '  -1.27%  '
$ws.Range("D48").Value = '1.68'
$ws.Range("E48").Value = '  -3.67%  '
$ws.Range("D49").Value = '88.15'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").Value = '0.0₆0108'
$ws.Range("E50").Value = '  -0.24%  '
$ws.Range("E51").Value = '  -3.13%  '
